# Update CDA Logical model for ST.r2b
#
# 1. Rename the "Include from ActRelationshipT" sheet to "Include #0"
# 2. Bump the Version / Date metadata values
# 3. Insert a new "Jurisdiction" row into the Metadata table (pushing
#    Description/Purpose/Copyright down by one row)
# 4. Append a new "Immutable" / "BooleanType[null]" row at the bottom of
#    the Metadata table

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1. Rename the Include tab -------------------------------------------
$ws2.Name = "Include #0"

# --- 2. Update Version and Date values on the Metadata sheet -------------
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- 3. Insert the new "Jurisdiction" row before the "Description" row ---
# (old row 11 "Description" / row 12 "Purpose" / row 13 "Copyright" all
# shift down by one)
$ws1.Rows.Item(11).Insert()

$ws1.Range("A11").Value = "Jurisdiction"
# Force a real (non-blank) empty string value, matching the source value
# that is an explicit empty string rather than an unset cell.
$ws1.Range("B11").Value = "'"

# Copy formatting (wrap text / border / fill) from the row that is now
# immediately below (the old "Description" row, now row 12) so the new
# row matches the sheet's existing style instead of Excel's insert default.
# (Applied after the values so PasteSpecial formats also clears any
# quote-prefix style flag picked up from the "'" assignment above.)
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Append the new "Immutable" row at the end of the table -----------
# Copy formatting from the row above (old "Immutable"/row 14, now row 14
# "Copyright") so the appended row matches the table's style.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A15").Value = "Immutable"
$ws1.Range("B15").Value = "BooleanType[null]"
